# Coposescu_Mihai_Octavian_prezentare.pptx
# "Advanced with pptx, modified figure inside synthesis"
#
# 1. Append 8 new "Title and Content" slides (Contextul proiectului ..
#    Bibliografie) after the existing "Continutul prezentarii" slide,
#    matching the table-of-contents order already listed there.
# 2. Refresh the cached datetimeFigureOut placeholder text (28.08.2024
#    -> 29.08.2024) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)

$CR = [char]13

function Add-ContentSlide($title, $bodyLines) {
    $s = $p.Slides.AddSlide($p.Slides.Count + 1, $titleAndContent)

    $s.Shapes.Item(1).TextFrame.TextRange.Text = $title

    if ($bodyLines -and $bodyLines.Length -gt 0) {
        $body = [string]::Join($CR, $bodyLines)
        $s.Shapes.Item(2).TextFrame.TextRange.Text = $body
    }
}

# Slide 3 - Contextul proiectului
Add-ContentSlide "Contextul proiectului" @(
    "Cum ne afecteaza calitatea aerului?",
    "",
    "90% din timp il petrecem in spatii inchise…",
    ""
) | Out-Null

# Slide 4 - Obiectivele proiectului
Add-ContentSlide "Obiectivele proiectului" @(
    "Masurarea parametrilor cheie de calitate a aerului.",
    "Citirea parametrilor in timp real.",
    "Citirea parametrilor istorici.",
    "Oferirea unei interfete intuitive si usor de utilizat."
) | Out-Null

# Slide 5 - Studiu bibliografic
Add-ContentSlide "Studiu bibliografic" @() | Out-Null

# Slide 6 - Solutia Aleasa
Add-ContentSlide "Solutia Aleasa" @() | Out-Null

# Slide 7 - Implementarea solutiei
Add-ContentSlide "Implementarea solutiei" @() | Out-Null

# Slide 8 - Teste si rezultate
Add-ContentSlide "Teste si rezultate" @() | Out-Null

# Slide 9 - Concluzii
Add-ContentSlide "Concluzii" @() | Out-Null

# Slide 10 - Bibliografie
Add-ContentSlide "Bibliografie" @() | Out-Null

# Refresh the cached "today" date text shown via the dt placeholder
# (type datetimeFigureOut) on the master and on every layout.
$newDate = "29.08.2024"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($li).Shapes
}
